$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1369.5454
$ws.Range("J17").Value = 1369.5454
$ws.Range("L17").Value = 4108.6362
$ws.Range("N17").Value = -4444.6362
$ws.Range("H38").Value = 186
$ws.Range("J38").Value = 444
$ws.Range("L38").Value = 1332
$ws.Range("N38").Value = -2076
$ws.Range("H41").Value = 499
$ws.Range("J41").Value = 499
$ws.Range("L41").Value = 499
$ws.Range("N41").Value = -1379
$ws.Range("H76").Value = 90917120
$ws.Range("J76").Value = 9899.5
$ws.Range("L76").Value = 9899.5
$ws.Range("N76").Value = -10529.5
$ws.Range("H79").Value = 90917120
$ws.Range("J79").Value = 9899.5
$ws.Range("L79").Value = 9899.5
$ws.Range("N79").Value = -12083.5
$ws.Range("H80").Value = 1057.4546
$ws.Range("I80").Value = 626
$ws.Range("J80").Value = 1417
$ws.Range("K80").Value = 1878
$ws.Range("L80").Value = 4251
$ws.Range("M80").Value = -880
$ws.Range("N80").Value = -6247
$ws.Range("H83").Value = 1057.4546
$ws.Range("I83").Value = 626
$ws.Range("J83").Value = 1417
$ws.Range("K83").Value = 5634
$ws.Range("L83").Value = 12753
$ws.Range("M83").Value = -642
$ws.Range("N83").Value = -22737
$ws.Range("H96").Value = 142858850
$ws.Range("J96").Value = 500002000
$ws.Range("L96").Value = 1500006000
$ws.Range("N96").Value = -1500008746
$ws.Range("H132").Value = 302582.1
$ws.Range("J132").Value = 10685.286
$ws.Range("L132").Value = 32055.858
$ws.Range("N132").Value = -37115.858
$ws.Range("H135").Value = 4259.8
$ws.Range("J135").Value = 9799.727999999999
$ws.Range("L135").Value = 88197.552
$ws.Range("N135").Value = -93267.552
$ws.Range("H138").Value = 3781.7778
$ws.Range("J138").Value = 5315.069
$ws.Range("L138").Value = 15945.207
$ws.Range("N138").Value = -26225.207
$ws.Range("H141").Value = 4274.25
$ws.Range("I141").Value = 3568.3845
$ws.Range("K141").Value = 10705.1535
$ws.Range("M141").Value = -5525.1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 25045000
$ws.Range("J9").Value = 89999
$ws.Range("L9").Value = 89999
$ws.Range("N9").Value = -90339
$ws.Range("H20").Value = 25045000
$ws.Range("J20").Value = 89999
$ws.Range("L20").Value = 89999
$ws.Range("N20").Value = -90539
$ws.Range("H102").Value = 5413.533
$ws.Range("I102").Value = 5443.0713
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 5443.0713
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -3821.0713
$ws.Range("N102").Value = -8244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12987892
$ws.Range("I20").Value = 14286570
$ws.Range("J20").Value = 1111
$ws.Range("K20").Value = 14286570
$ws.Range("L20").Value = 1111
$ws.Range("M20").Value = -14286323
$ws.Range("N20").Value = -1605
$ws.Range("H86").Value = 22535.143
$ws.Range("I86").Value = 50450
$ws.Range("J86").Value = 11369.2
$ws.Range("K86").Value = 50450
$ws.Range("L86").Value = 11369.2
$ws.Range("M86").Value = -49327
$ws.Range("N86").Value = -13615.2
$ws.Range("H89").Value = 22535.143
$ws.Range("I89").Value = 50450
$ws.Range("J89").Value = 11369.2
$ws.Range("K89").Value = 252250
$ws.Range("L89").Value = 56846
$ws.Range("M89").Value = -246634
$ws.Range("N89").Value = -68078
$ws.Range("H99").Value = 6943.8184
$ws.Range("I99").Value = 1655.6
$ws.Range("K99").Value = 1655.6
$ws.Range("M99").Value = -157.5999999999999
$ws.Range("H105").Value = 2916.5
$ws.Range("J105").Value = 2800.5806
$ws.Range("L105").Value = 2800.5806
$ws.Range("N105").Value = -6294.580599999999
$ws.Range("H134").Value = 1258460.9
$ws.Range("I134").Value = 1447492.5
$ws.Range("J134").Value = 10852.6
$ws.Range("K134").Value = 4342477.5
$ws.Range("L134").Value = 32557.8
$ws.Range("M134").Value = -4339942.5
$ws.Range("N134").Value = -37627.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4513.561
$ws.Range("I132").Value = 3561.8286
$ws.Range("K132").Value = 10685.4858
$ws.Range("M132").Value = -8155.485799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5637.5
$ws.Range("I34").Value = 50
$ws.Range("J34").Value = 7500
$ws.Range("K34").Value = 150
$ws.Range("L34").Value = 22500
$ws.Range("M34").Value = -66
$ws.Range("N34").Value = -22668
$ws.Range("H41").Value = 6666
$ws.Range("H50").Value = 3600.8
$ws.Range("I50").Value = 3251
$ws.Range("J50").Value = 5000
$ws.Range("K50").Value = 9753
$ws.Range("L50").Value = 15000
$ws.Range("M50").Value = -9272
$ws.Range("N50").Value = -15962
$ws.Range("H53").Value = 3600.8
$ws.Range("I53").Value = 3251
$ws.Range("J53").Value = 5000
$ws.Range("K53").Value = 9753
$ws.Range("L53").Value = 15000
$ws.Range("M53").Value = -9272
$ws.Range("N53").Value = -15962
$ws.Range("H109").Value = 10730.214
$ws.Range("I109").Value = 7358.1113
$ws.Range("J109").Value = 16800
$ws.Range("K109").Value = 22074.3339
$ws.Range("L109").Value = 50400
$ws.Range("M109").Value = -21034.3339
$ws.Range("N109").Value = -52480
$ws.Range("H119").Value = 2526.2856
$ws.Range("I119").Value = 2526.2856
$ws.Range("K119").Value = 7578.8568
$ws.Range("M119").Value = -2740.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1469.6666
$ws.Range("J107").Value = 2651.5
$ws.Range("L107").Value = 2651.5
$ws.Range("N107").Value = -6491.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 8001.4
$ws.Range("I21").Value = 2000
$ws.Range("J21").Value = 9501.75
$ws.Range("K21").Value = 2000
$ws.Range("L21").Value = 9501.75
$ws.Range("M21").Value = -1826
$ws.Range("N21").Value = -9849.75
$ws.Range("H82").Value = 4415.067
$ws.Range("I82").Value = 2300.6667
$ws.Range("J82").Value = 7586.6665
$ws.Range("K82").Value = 2300.6667
$ws.Range("L82").Value = 7586.6665
$ws.Range("M82").Value = -1939.6667
$ws.Range("N82").Value = -8308.666499999999
$ws.Range("H85").Value = 4415.067
$ws.Range("I85").Value = 2300.6667
$ws.Range("J85").Value = 7586.6665
$ws.Range("K85").Value = 2300.6667
$ws.Range("L85").Value = 7586.6665
$ws.Range("M85").Value = -1052.6667
$ws.Range("N85").Value = -10082.6665
$ws.Range("H106").Value = 17283.5
$ws.Range("J106").Value = 17283.5
$ws.Range("L106").Value = 17283.5
$ws.Range("N106").Value = -19807.5
$ws.Range("H132").Value = 2272.1428
$ws.Range("I132").Value = 2272.1428
$ws.Range("K132").Value = 6816.428400000001
$ws.Range("M132").Value = -4286.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16999.572
$ws.Range("I62").Value = 17249.666
$ws.Range("J62").Value = 16812
$ws.Range("K62").Value = 17249.666
$ws.Range("L62").Value = 16812
$ws.Range("M62").Value = -16625.666
$ws.Range("N62").Value = -18060
$ws.Range("H65").Value = 16999.572
$ws.Range("I65").Value = 17249.666
$ws.Range("J65").Value = 16812
$ws.Range("K65").Value = 86248.33
$ws.Range("L65").Value = 84060
$ws.Range("M65").Value = -83128.33
$ws.Range("N65").Value = -90300
$ws.Range("H81").Value = 1575
$ws.Range("I81").Value = 1350
$ws.Range("J81").Value = 2250
$ws.Range("K81").Value = 2700
$ws.Range("L81").Value = 4500
$ws.Range("M81").Value = -1639
$ws.Range("N81").Value = -6622
$ws.Range("H84").Value = 1575
$ws.Range("I84").Value = 1350
$ws.Range("J84").Value = 2250
$ws.Range("K84").Value = 13500
$ws.Range("L84").Value = 22500
$ws.Range("M84").Value = -8196
$ws.Range("N84").Value = -33108
$ws.Range("H104").Value = 19376.9
$ws.Range("J104").Value = 19376.9
$ws.Range("L104").Value = 19376.9
$ws.Range("N104").Value = -26364.9
